# Change data for test purpose
$wb = $excel.ActiveWorkbook

# --- DRONE2 sheet ---
$ws2 = $wb.Worksheets.Item("DRONE2")
$ws2.Activate()
$ws2.Range("A2").Value = 112.3
$ws2.Range("A2").Select() | Out-Null

# --- DRONE3 sheet ---
$ws3 = $wb.Worksheets.Item("DRONE3")
$ws3.Activate()
$ws3.Range("B2").Value = 0.28499999999999998
# 5.1666666666666667 (character width) renders/stores as exactly 6 (stored width units)
$ws3.Columns("B:B").ColumnWidth = 5.1666666666666667
$ws3.Cells.Select() | Out-Null
